# update VOM, now include fuel costs
$wb = $excel.ActiveWorkbook

# Work on the FOAK sheet, which holds the VOM ($/MWh-e) column (O)
$ws = $wb.Worksheets.Item("FOAK")

# Update VOM values to include fuel costs
$ws.Range("O2").Formula = "=0.75+9"
$ws.Range("O3").Value = 13
$ws.Range("O4").Value = 13
$ws.Range("O5").Value = 11.5

# Make FOAK the active sheet, scroll/select as the author left it
$ws.Activate()
$ws.Range("O7").Select()
